$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix capitalization typos ("intoduction" -> "Introduction") in the task names
# for the three most recent lecture entries.
$ws.Cells.Item(13, 6).Value = "CS Introduction Lecture 12"
$ws.Cells.Item(14, 6).Value = "CS Introduction Lecture 13"
$ws.Cells.Item(15, 6).Value = "CS Introduction Lecture 14"

# Add a new log entry for Lecture 15 in row 16, matching the date-formatted
# style already used in column A.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(16, 1).Value = 45814
$ws.Cells.Item(16, 2).Value = 14
$ws.Cells.Item(16, 3).Value = 45
$ws.Cells.Item(16, 4).Value = 15
$ws.Cells.Item(16, 5).Value = 30
$ws.Cells.Item(16, 6).Value = "CS Introduction Lecture 15"

# Move the active selection to reflect where the user left off after
# entering the new row.
$ws.Range("E17").Select()
